# edit.ps1 - applies the evolang_XV_template.docx content edit described by the
# commit "Update page limit in templates and homogenize instructions".
#
# Real, observable changes made to the document:
#   1. The "Abstract" style paragraph in the instructions is reworded:
#        - both the short-submission and long-submission page-limit
#          parentheticals now also exclude "acknowledgements", not just
#          "references";
#        - a stray " ," / extra space after the long-submission
#          parenthetical is cleaned up;
#        - "number - indicate" -> "number \u2013 indicate" (en dash);
#        - "title \u2013 acknowledgements ... in a separate section" becomes
#          "title. Acknowledgements ... in the separate 'Acknowledgements'
#          section", explicitly naming the Acknowledgements section.
#   2. A stale cached pagination hint (<w:lastRenderedPageBreak/>) in the
#      Acknowledgements/Appendices paragraph is cleared.
#   3. The section's page size is made to carry an explicit portrait
#      orientation flag.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Abstract paragraph rewording
# ---------------------------------------------------------------------

$enDash = [char]0x2013
$lsquo  = [char]0x2018
$rsquo  = [char]0x2019

# 1a. short-submission parenthetical gains "and acknowledgements"
$null = $d.Content.Find.Execute(
    "(2 pages excluding references),", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(2 pages excluding references and acknowledgements),", 2)

# 1b. long-submission parenthetical gains "and acknowledgements" and the
#     stray space/comma spacing around it is tidied up
$null = $d.Content.Find.Execute(
    "(6 pages excluding references) ,", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(6 pages excluding references and acknowledgements),", 2)

# 1c. hyphen -> en dash before "indicate"
$null = $d.Content.Find.Execute(
    "superscript number - indicate", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "superscript number $enDash indicate", 2)

# 1d. point to the named Acknowledgements section
$oldTail = "abstract or title $enDash acknowledgements for funding bodies etc. are to be placed in a separate section"
$newTail = "abstract or title. Acknowledgements for funding bodies etc. are to be placed in the separate ${lsquo}Acknowledgements$rsquo section"
$null = $d.Content.Find.Execute(
    $oldTail, $true, $false, $false, $false, $false,
    $true, 1, $false, $newTail, 2)

# ---------------------------------------------------------------------
# 2. Drop the stale lastRenderedPageBreak cached in the Acknowledgements /
#    Appendices paragraph (harmless pagination cache, recomputed on layout).
# ---------------------------------------------------------------------

$null = $d.Content.Find.Execute(
    "hosted online along with the final proceedings", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "hosted online along with the final proceedings", 2)

# ---------------------------------------------------------------------
# 3. Make the page orientation explicit (portrait) on the section.
# ---------------------------------------------------------------------

$d.PageSetup.Orientation = 0   # wdOrientPortrait
